# Update countries & provincias Spain
#
# The source "paises" (countries) list behind this sheet was re-pulled, which
# (a) refreshed the case/death counters for a handful of countries and
# (b) changed the sort order (rows are ranked by column B, "Casos totales",
#     descending), so several countries swapped rows with their neighbours:
#       - Turquia overtook Pakistan                         (rows 16-17)
#       - Paraguay overtook Guinea-Bisau / Eslovenia         (rows 117-119)
#       - Namibia overtook Monaco / Aruba / Barbados / Botsuana (rows 179-183)
#       - Fiyi overtook Dominica                             (rows 202-203)
#       - Groenlandia overtook Islas Malvinas                (rows 208-209)
#       - Seychelles overtook Montserrat                     (rows 211-212)
# The footer timestamp was also bumped to the new refresh time.
#
# Column layout (row 3 header): A Pais | B Casos totales | C Nuevos casos |
#   D Casos activos | E Recuperados | F Casos criticos | G Muertes hoy | H Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer "last updated" timestamp (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Junio de 2020 a las 20:36"

# Estados Unidos (row 4) - refreshed counters
$ws.Cells.Item(4, 2).Value = 2483398
$ws.Cells.Item(4, 3).Value = 20844
$ws.Cells.Item(4, 4).Value = 1042266
$ws.Cells.Item(4, 5).Value = 1316544
$ws.Cells.Item(4, 7).Value = 307
$ws.Cells.Item(4, 8).Value = 124588

# India (row 7) - refreshed counters
$ws.Cells.Item(7, 2).Value = 490892
$ws.Cells.Item(7, 3).Value = 17907
$ws.Cells.Item(7, 4).Value = 285664
$ws.Cells.Item(7, 5).Value = 189920
$ws.Cells.Item(7, 7).Value = 401
$ws.Cells.Item(7, 8).Value = 15308

# Espana (row 9) - refreshed counters
$ws.Cells.Item(9, 2).Value = 294566
$ws.Cells.Item(9, 3).Value = 400
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 8).Value = 28330

# Turquia now ranks above Pakistan -> row 16 becomes Turquia
$ws.Cells.Item(16, 1).Value = "Turquia"
$ws.Cells.Item(16, 2).Value = 193115
$ws.Cells.Item(16, 3).Value = 1458
$ws.Cells.Item(16, 4).Value = 165706
$ws.Cells.Item(16, 5).Value = 22363
$ws.Cells.Item(16, 7).Value = 21
$ws.Cells.Item(16, 8).Value = 5046

# ...and Pakistan drops to row 17
$ws.Cells.Item(17, 1).Value = "Pakistan"
$ws.Cells.Item(17, 2).Value = 192970
$ws.Cells.Item(17, 3).Value = 4044
$ws.Cells.Item(17, 4).Value = 81307
$ws.Cells.Item(17, 5).Value = 107760
$ws.Cells.Item(17, 7).Value = 148
$ws.Cells.Item(17, 8).Value = 3903

# Chequia (row 50) - refreshed counters
$ws.Cells.Item(50, 5).Value = 5522
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 8).Value = 71

# Costa Rica (row 67) - refreshed counters
$ws.Cells.Item(67, 2).Value = 11338
$ws.Cells.Item(67, 3).Value = 431
$ws.Cells.Item(67, 4).Value = 8500
$ws.Cells.Item(67, 5).Value = 2621
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = 217

# Estonia (row 104) - refreshed counters
$ws.Cells.Item(104, 2).Value = 2277
$ws.Cells.Item(104, 3).Value = 16
$ws.Cells.Item(104, 4).Value = 1848
$ws.Cells.Item(104, 5).Value = 421

# Paraguay now ranks above Guinea-Bisau and Eslovenia -> row 117 becomes Paraguay
$ws.Cells.Item(117, 1).Value = "Paraguay"
$ws.Cells.Item(117, 2).Value = 1569
$ws.Cells.Item(117, 3).Value = 41
$ws.Cells.Item(117, 4).Value = 976
$ws.Cells.Item(117, 5).Value = 580
$ws.Cells.Item(117, 8).Value = 13

# ...Guinea-Bisau drops to row 118
$ws.Cells.Item(118, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(118, 2).Value = 1556
$ws.Cells.Item(118, 3).Value = 0
$ws.Cells.Item(118, 4).Value = 191
$ws.Cells.Item(118, 5).Value = 1346
$ws.Cells.Item(118, 8).Value = 19

# ...and Eslovenia drops to row 119
$ws.Cells.Item(119, 1).Value = "Eslovenia"
$ws.Cells.Item(119, 2).Value = 1547
$ws.Cells.Item(119, 3).Value = 6
$ws.Cells.Item(119, 4).Value = 1376
$ws.Cells.Item(119, 5).Value = 62
$ws.Cells.Item(119, 8).Value = 109

# Nueva Zelanda (row 122) - refreshed counters
$ws.Cells.Item(122, 2).Value = 1382
$ws.Cells.Item(122, 3).Value = 54
$ws.Cells.Item(122, 4).Value = 446
$ws.Cells.Item(122, 5).Value = 933

# Letonia (row 134) - refreshed counters
$ws.Cells.Item(134, 2).Value = 960
$ws.Cells.Item(134, 3).Value = 19
$ws.Cells.Item(134, 4).Value = 260
$ws.Cells.Item(134, 5).Value = 688
$ws.Cells.Item(134, 7).Value = 1
$ws.Cells.Item(134, 8).Value = 12

# Namibia now ranks above Monaco, Aruba, Barbados and Botsuana -> row 179 becomes Namibia
$ws.Cells.Item(179, 1).Value = "Namibia"
$ws.Cells.Item(179, 3).Value = 26
$ws.Cells.Item(179, 4).Value = 21
$ws.Cells.Item(179, 5).Value = 81
$ws.Cells.Item(179, 8).Value = 0

# ...Monaco drops to row 180
$ws.Cells.Item(180, 1).Value = "Monaco"
$ws.Cells.Item(180, 2).Value = 102
$ws.Cells.Item(180, 4).Value = 95
$ws.Cells.Item(180, 5).Value = 3
$ws.Cells.Item(180, 8).Value = 4

# ...Aruba drops to row 181
$ws.Cells.Item(181, 1).Value = "Aruba"
$ws.Cells.Item(181, 2).Value = 101
$ws.Cells.Item(181, 4).Value = 98
$ws.Cells.Item(181, 5).Value = 0
$ws.Cells.Item(181, 8).Value = 3

# ...Barbados drops to row 182
$ws.Cells.Item(182, 1).Value = "Barbados"
$ws.Cells.Item(182, 2).Value = 97
$ws.Cells.Item(182, 4).Value = 85
$ws.Cells.Item(182, 5).Value = 5
$ws.Cells.Item(182, 8).Value = 7

# ...and Botsuana drops to row 183
$ws.Cells.Item(183, 1).Value = "Botsuana"
$ws.Cells.Item(183, 2).Value = 92
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 25
$ws.Cells.Item(183, 5).Value = 66
$ws.Cells.Item(183, 8).Value = 1

# Fiyi now ranks above Dominica -> row 202/203 swap (counters unchanged for these two)
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"

# Groenlandia now ranks above Islas Malvinas -> row 208/209 swap (counters unchanged)
$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"

# Seychelles now ranks above Montserrat -> row 211 becomes Seychelles
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0

# ...and Montserrat drops to row 212
$ws.Cells.Item(212, 1).Value = "Montserrat"
$ws.Cells.Item(212, 4).Value = 10
$ws.Cells.Item(212, 8).Value = 1
